# "Generate Report for Handoff" - the localization status report moved from
# "In Translation" to "Ready for handoff", so refresh the status + the
# handoff timestamps on all three sheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-27 04:37:51"

# Status text got longer ("In Translation" -> "Ready for handoff"), so the
# zh-cn / de-de status columns widen to fit it.
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333

# --- zh-cn sheet -------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Ready for handoff"
$wsZh.Range("H2").Value = "2016-08-27 04:37:46"
$wsZh.Columns.Item(3).ColumnWidth = 16.333333

# --- de-de sheet -------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Ready for handoff"
$wsDe.Range("H2").Value = "2016-08-27 04:37:51"
$wsDe.Columns.Item(3).ColumnWidth = 16.333333
